# Reorder columns of the Learning Objectives sheet:
#   old: A=Topic, B=LO Code, C=LO Description, D=Subject, E=Grade
#   new: A=Grade, B=Subject, C=Topic, D=LO Code, E=LO Description
# and rename headers to lowercase/underscored form.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count

$srcRange = $ws.Range($used.Address())
$data = $srcRange.Value()

$newData = New-Object 'object[,]' $rowCount,5

for ($r = 1; $r -le $rowCount; $r++) {
    $newData[$r-1,0] = $data[$r,5]
    $newData[$r-1,1] = $data[$r,4]
    $newData[$r-1,2] = $data[$r,1]
    $newData[$r-1,3] = $data[$r,2]
    $newData[$r-1,4] = $data[$r,3]
}

$newData[0,0] = "grade"
$newData[0,1] = "subject"
$newData[0,2] = "topic"
$newData[0,3] = "lo_code"
$newData[0,4] = "lo_description"

$destRange = $ws.Range($used.Address())

# Force the LO-code column (new column D, e.g. "11.1.10") to Text format so
# Excel's automatic data-type detection doesn't reinterpret values such as
# "11.1.10" as a date when they are written back through .Value.
$codeColRange = $ws.Range("D1:D$rowCount")
$codeColRange.NumberFormat = "@"

$destRange.Value = $newData

$ws.Range("A1:E1").Select()
